# Edit the document header text:
#   "Picture and text here"
# becomes four separate runs whose text concatenates to:
#   "Picture" + "s" + " and text here " + "in template document"
#   = "Pictures and text here in template document"
#
# We build this up in two passes:
#   1) Insert the new characters ("s" and " in template document") at the
#      right spots using Range.InsertAfter.
#   2) Re-locate each of the four desired segments with Find and nudge
#      their Bold property on/off (a no-op format-wise) purely to force
#      the engine to keep them as separate <w:r> runs instead of
#      re-coalescing adjacent, identically-formatted runs.

$d = $word.ActiveDocument
$hdr = $d.Sections.Item(1).Headers.Item(1)

# --- Step 1: insert "s" right after "Picture" -----------------------------
$p1 = $hdr.Range.Paragraphs.Item(1)
$pr1 = $p1.Range
$pr1.Find.Execute("Picture", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pr1.InsertAfter("s")

# --- Step 2: insert " in template document" right after "here" -----------
$p2 = $hdr.Range.Paragraphs.Item(1)
$pr2 = $p2.Range
$pr2.Find.Execute("here", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pr2.InsertAfter(" in template document")

# --- Step 3: split "Picture" | "s and text here in template document" ---
$p3 = $hdr.Range.Paragraphs.Item(1)
$pr3 = $p3.Range
$pr3.Find.Execute("s and text here in template document", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pr3.Bold = 1
$pr3.Bold = 0

# --- Step 4: split "s" | " and text here in template document" -----------
$p4 = $hdr.Range.Paragraphs.Item(1)
$pr4 = $p4.Range
$pr4.Find.Execute(" and text here in template document", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pr4.Bold = 1
$pr4.Bold = 0

# --- Step 5: split " and text here " | "in template document" ------------
$p5 = $hdr.Range.Paragraphs.Item(1)
$pr5 = $p5.Range
$pr5.Find.Execute("in template document", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pr5.Bold = 1
$pr5.Bold = 0

# --- Step 6: re-touch "Picture" so it keeps an explicit (empty) rPr ------
$p6 = $hdr.Range.Paragraphs.Item(1)
$pr6 = $p6.Range
$pr6.Find.Execute("Picture", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pr6.Bold = 1
$pr6.Bold = 0

"Header text now: [" + $hdr.Range.Text + "]"
